$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Whisp outputs")

# Insert a new row above row 166 - this shifts existing rows 166..195 down to 167..196
$ws.Rows.Item(166).Insert()
$ws.Rows.Item(166).RowHeight = 18.75

# Populate the newly inserted row 166 with the new dataset entry
$ws.Cells.Item(166, 1).Value = "DIST_after_2020"
$ws.Cells.Item(166, 2).Value = "numeric"
$ws.Cells.Item(166, 3).Value = "ha / %"
$ws.Cells.Item(166, 4).Value = "Area of Tree cover loss"
$ws.Cells.Item(166, 5).Value = "2023-2026"
$ws.Cells.Item(166, 6).Value = "Pickens 2025"
$ws.Cells.Item(166, 7).Value = $null
$ws.Cells.Item(166, 8).Value = 1
$ws.Cells.Item(166, 9).Value = "disturbance_after"
$ws.Cells.Item(166, 10).Value = 1
$ws.Cells.Item(166, 11).Value = "disturbance_after"
$ws.Cells.Item(166, 12).Value = 0
$ws.Cells.Item(166, 13).Value = $null
$ws.Cells.Item(166, 14).Value = $null

# Reflect the saved view/selection state from the author's edit
$ws.Activate() | Out-Null
$ws.Range("A166").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 159 | Out-Null
